$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.080686369548288
$ws.Range("D2").Value = 1.081314231955256
$ws.Range("E2").Value = 1.083625367113356
$ws.Range("F2").Value = 1.093730810571756
$ws.Range("I2").Value = 1.067001903477116
$ws.Range("J2").Value = 1.085564674350209
$ws.Range("K2").Value = 1.083985378331075
$ws.Range("L2").Value = 1.086290496804683
$ws.Range("M2").Value = 1.096369954386104
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.081900152883505
$ws.Range("D3").Value = 1.082279443728756
$ws.Range("E3").Value = 1.084690541199117
$ws.Range("F3").Value = 1.094823840609915
$ws.Range("I3").Value = 1.06744291583022
$ws.Range("J3").Value = 1.08643865572687
$ws.Range("K3").Value = 1.084768655942776
$ws.Range("L3").Value = 1.087173930016216
$ws.Range("M3").Value = 1.097283056065762
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.082685458684607
$ws.Range("D4").Value = 1.082903842965748
$ws.Range("E4").Value = 1.085379832326296
$ws.Range("F4").Value = 1.095531235005786
$ws.Range("I4").Value = 1.067726974172488
$ws.Range("J4").Value = 1.087003501696314
$ws.Range("K4").Value = 1.085274706203559
$ws.Range("L4").Value = 1.087745016320181
$ws.Range("M4").Value = 1.09787342173235
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.083015580419484
$ws.Range("D5").Value = 1.083166303299428
$ws.Range("E5").Value = 1.085669623548776
$ws.Range("F5").Value = 1.095828655638713
$ws.Range("I5").Value = 1.067846080079299
$ws.Range("J5").Value = 1.087240801340487
$ws.Range("K5").Value = 1.085487262998473
$ws.Range("L5").Value = 1.087984969089828
$ws.Range("M5").Value = 1.09812149911814
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.083071008133467
$ws.Range("D6").Value = 1.083210369384688
$ws.Range("E6").Value = 1.085718281567111
$ws.Range("F6").Value = 1.095878595770215
$ws.Range("I6").Value = 1.067866060190167
$ws.Range("J6").Value = 1.087280635542254
$ws.Range("K6").Value = 1.08552294128843
$ws.Range("L6").Value = 1.088025250501551
$ws.Range("M6").Value = 1.098163145835498
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.082689869867342
$ws.Range("D7").Value = 1.082907350118672
$ws.Range("E7").Value = 1.085383704477893
$ws.Range("F7").Value = 1.09553520902811
$ws.Range("I7").Value = 1.067728566897663
$ws.Range("J7").Value = 1.087006673139734
$ws.Range("K7").Value = 1.085277547130777
$ws.Range("L7").Value = 1.087748223098242
$ws.Range("M7").Value = 1.097876736995063
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.081096593495041
$ws.Range("D8").Value = 1.081640462703723
$ws.Range("E8").Value = 1.083985337293443
$ws.Range("F8").Value = 1.094100178512695
$ws.Range("I8").Value = 1.067151216395424
$ws.Range("J8").Value = 1.085860181410224
$ws.Range("K8").Value = 1.084250253129277
$ws.Range("L8").Value = 1.086589172355007
$ws.Range("M8").Value = 1.096678639562402
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.078288256736938
$ws.Range("D9").Value = 1.079406814918398
$ws.Range("E9").Value = 1.081521603372533
$ws.Range("F9").Value = 1.091572442080404
$ws.Range("I9").Value = 1.066123826274
$ws.Range("J9").Value = 1.083834683858763
$ws.Range("K9").Value = 1.082434005563908
$ws.Range("L9").Value = 1.084542500478427
$ws.Range("M9").Value = 1.09456379083628
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.076415412255022
$ws.Range("D10").Value = 1.077916851078963
$ws.Range("E10").Value = 1.079879313134068
$ws.Range("F10").Value = 1.089887890111293
$ws.Range("I10").Value = 1.065432128266146
$ws.Range("J10").Value = 1.082480780453501
$ws.Range("K10").Value = 1.081219080078008
$ws.Range("L10").Value = 1.08317512902092
$ws.Range("M10").Value = 1.093151392348075
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.075604280900415
$ws.Range("D11").Value = 1.077271465728615
$ws.Range("E11").Value = 1.079168218982718
$ws.Range("F11").Value = 1.089158592844757
$ws.Range("I11").Value = 1.065131001947356
$ws.Range("J11").Value = 1.081893665542158
$ws.Range("K11").Value = 1.080692022600123
$ws.Range("L11").Value = 1.082582337180067
$ws.Range("M11").Value = 1.092539204261199
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.075302961965503
$ws.Range("D12").Value = 1.077031706740822
$ws.Range("E12").Value = 1.078904090117167
$ws.Range("F12").Value = 1.088887717047013
$ws.Range("I12").Value = 1.065018906653416
$ws.Range("J12").Value = 1.081675453662157
$ws.Range("K12").Value = 1.080496100757733
$ws.Range("L12").Value = 1.08236203999572
$ws.Range("M12").Value = 1.092311717611418
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.075367597227631
$ws.Range("D13").Value = 1.077083137416699
$ws.Range("E13").Value = 1.078960746529893
$ws.Range("F13").Value = 1.088945820042462
$ws.Range("I13").Value = 1.065042962507965
$ws.Range("J13").Value = 1.081722266809984
$ws.Range("K13").Value = 1.080538133426669
$ws.Range("L13").Value = 1.082409299389721
$ws.Range("M13").Value = 1.092360518483538
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.075579374378903
$ws.Range("D14").Value = 1.077251647877422
$ws.Range("E14").Value = 1.079146385950102
$ws.Range("F14").Value = 1.08913620181241
$ws.Range("I14").Value = 1.065121741088266
$ws.Range("J14").Value = 1.081875630764608
$ws.Range("K14").Value = 1.080675830689404
$ws.Range("L14").Value = 1.082564129553497
$ws.Range("M14").Value = 1.092520402041681
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.075709853417458
$ws.Range("D15").Value = 1.077355468183827
$ws.Range("E15").Value = 1.079260764910173
$ws.Range("F15").Value = 1.089253504614627
$ws.Range("I15").Value = 1.065170246874864
$ws.Range("J15").Value = 1.081970105929991
$ws.Range("K15").Value = 1.080760650713998
$ws.Range("L15").Value = 1.082659511197132
$ws.Range("M15").Value = 1.092618899258527
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.076469240465384
$ws.Range("D16").Value = 1.077959678442115
$ws.Range("E16").Value = 1.079926506614845
$ws.Range("F16").Value = 1.089936293707443
$ws.Range("I16").Value = 1.065452078909569
$ws.Range("J16").Value = 1.082519726994061
$ws.Range("K16").Value = 1.081254038273876
$ws.Range("L16").Value = 1.083214455591061
$ws.Range("M16").Value = 1.093192008308267
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.076945535014757
$ws.Range("D17").Value = 1.078338623721978
$ws.Range("E17").Value = 1.080344115616411
$ws.Range("F17").Value = 1.090364621932187
$ws.Range("I17").Value = 1.065628431411615
$ws.Range("J17").Value = 1.082864257131751
$ws.Range("K17").Value = 1.081563262386784
$ws.Range("L17").Value = 1.083562366434176
$ws.Range("M17").Value = 1.093551340354094
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.077223332590492
$ws.Range("D18").Value = 1.078559634812968
$ws.Range("E18").Value = 1.08058770269156
$ws.Range("F18").Value = 1.090614470657842
$ws.Range("I18").Value = 1.065731138919414
$ws.Range("J18").Value = 1.083065132274109
$ws.Range("K18").Value = 1.081743532441508
$ws.Range("L18").Value = 1.083765228416611
$ws.Range("M18").Value = 1.093760873850521
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.077318051565926
$ws.Range("D19").Value = 1.078634990301644
$ws.Range("E19").Value = 1.08067076009929
$ws.Range("F19").Value = 1.090699664704981
$ws.Range("I19").Value = 1.065766133099684
$ws.Range("J19").Value = 1.083133611386653
$ws.Range("K19").Value = 1.081804983764316
$ws.Range("L19").Value = 1.083834387519954
$ws.Range("M19").Value = 1.093832309415538
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.076894434887316
$ws.Range("D20").Value = 1.07829796868714
$ws.Range("E20").Value = 1.080299309839842
$ws.Range("F20").Value = 1.090318665121411
$ws.Range("I20").Value = 1.065609526589977
$ws.Range("J20").Value = 1.082827300949388
$ws.Range("K20").Value = 1.081530095413514
$ws.Range("L20").Value = 1.083525045988942
$ws.Range("M20").Value = 1.09351279351635
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.075517012092698
$ws.Range("D21").Value = 1.077202026696253
$ws.Range("E21").Value = 1.079091719689271
$ws.Range("F21").Value = 1.089080138646942
$ws.Range("I21").Value = 1.065098549489318
$ws.Range("J21").Value = 1.08183047253541
$ws.Range("K21").Value = 1.080635286404862
$ws.Range("L21").Value = 1.082518538916599
$ws.Range("M21").Value = 1.092473322887283
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.074650805057516
$ws.Range("D22").Value = 1.076512766537887
$ws.Range("E22").Value = 1.078332477593226
$ws.Range("F22").Value = 1.088301529645893
$ws.Range("I22").Value = 1.064775868326183
$ws.Range("J22").Value = 1.081202966927916
$ws.Range("K22").Value = 1.080071820592794
$ws.Range("L22").Value = 1.081885083155689
$ws.Range("M22").Value = 1.091819229376716
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.075110014258161
$ws.Range("D23").Value = 1.076878175452934
$ws.Range("E23").Value = 1.078734964876799
$ws.Range("F23").Value = 1.088714275724995
$ws.Range("I23").Value = 1.064947061564191
$ws.Range("J23").Value = 1.081535691959732
$ws.Range("K23").Value = 1.080370606694717
$ws.Range("L23").Value = 1.082220949627523
$ws.Range("M23").Value = 1.092166027983997
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.076917524887247
$ws.Range("D24").Value = 1.078316339012849
$ws.Range("E24").Value = 1.080319555632598
$ws.Range("F24").Value = 1.090339430988218
$ws.Range("I24").Value = 1.065618069346679
$ws.Range("J24").Value = 1.082844000114949
$ws.Range("K24").Value = 1.081545082434729
$ws.Range("L24").Value = 1.083541909703099
$ws.Range("M24").Value = 1.093530211354463
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.079014381391029
$ws.Range("D25").Value = 1.079984416128499
$ws.Range("E25").Value = 1.082158499029559
$ws.Range("F25").Value = 1.092225812340163
$ws.Range("I25").Value = 1.066390622446641
$ws.Range("J25").Value = 1.084358949137674
$ws.Range("K25").Value = 1.082904267470118
$ws.Range("L25").Value = 1.085072125968345
$ws.Range("M25").Value = 1.095110967029838
